$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.339.68"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "2.510.26"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.56"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.33"
$ws.Range("E6").Value = "  -4.34%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E8").Value = "  -1.72%  "
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("E12").Value = "  -3.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.354"
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").Value = "2.964.72"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.37"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "59.249.02"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "2.512.19"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.08"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.30"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.35"
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("E25").Value = "  -3.33%  "
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.80"
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("D30").Value = "0.0₃0778"
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("E31").Value = "  -1.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.59"
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.44"
$ws.Range("E34").Value = "  -3.27%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.12"
$ws.Range("E35").Value = "  -7.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.49"
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.22"
$ws.Range("E37").Value = "  -4.19%  "
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.92"
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("E41").Value = "  -3.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.22"
$ws.Range("E42").Value = "  -7.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "280.99"
$ws.Range("E43").Value = "  -5.49%  "
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.596"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "124.94"
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0936"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0224"
$ws.Range("E50").Value = "  -1.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.83"
$ws.Range("E51").Value = "  -2.81%  "
